# Auto-generated edit script: update leve-profit value cells per scheduled data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 805.55554
$ws.Range("I4").Value = 792.8570999999999
$ws.Range("J4").Value = 850
$ws.Range("K4").Value = 792.8570999999999
$ws.Range("L4").Value = 850
$ws.Range("M4").Value = -678.8570999999999
$ws.Range("N4").Value = -1078

$ws.Range("H111").Value = 7649.8887
$ws.Range("I111").Value = 4856.125
$ws.Range("K111").Value = 14568.375
$ws.Range("M111").Value = -11501.375

$ws.Range("H112").Value = 2587.4583
$ws.Range("I112").Value = 841.6667
$ws.Range("J112").Value = 3169.389
$ws.Range("K112").Value = 2525.0001
$ws.Range("L112").Value = 9508.167000000001
$ws.Range("M112").Value = -1417.0001
$ws.Range("N112").Value = -11724.167

$ws.Range("H113").Value = 4124.7617
$ws.Range("I113").Value = 2851.25
$ws.Range("J113").Value = 4424.4116
$ws.Range("K113").Value = 2851.25
$ws.Range("L113").Value = 4424.4116
$ws.Range("M113").Value = 402.75
$ws.Range("N113").Value = -10932.4116

$ws.Range("H137").Value = 13515516
$ws.Range("I137").Value = 19232774
$ws.Range("K137").Value = 57698322
$ws.Range("M137").Value = -57695772

$ws.Range("H138").Value = 1987.1455
$ws.Range("I138").Value = 969.5
$ws.Range("J138").Value = 2899.5173
$ws.Range("K138").Value = 2908.5
$ws.Range("L138").Value = 8698.5519
$ws.Range("M138").Value = 2231.5
$ws.Range("N138").Value = -18978.5519

$ws.Range("H141").Value = 13132
$ws.Range("I141").Value = 3747.9
$ws.Range("J141").Value = 60052.5
$ws.Range("K141").Value = 11243.7
$ws.Range("L141").Value = 180157.5
$ws.Range("M141").Value = -6063.700000000001
$ws.Range("N141").Value = -190517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2217.76
$ws.Range("I2").Value = 1349
$ws.Range("J2").Value = 4968.8335
$ws.Range("K2").Value = 1349
$ws.Range("L2").Value = 4968.8335
$ws.Range("M2").Value = -1236
$ws.Range("N2").Value = -5194.8335

$ws.Range("H32").Value = 4049.0598
$ws.Range("I32").Value = 3965.4219
$ws.Range("J32").Value = 5833.3335
$ws.Range("K32").Value = 3965.4219
$ws.Range("L32").Value = 5833.3335
$ws.Range("M32").Value = -3678.4219
$ws.Range("N32").Value = -6407.3335

$ws.Range("H45").Value = 3585.35
$ws.Range("I45").Value = 2482.8572
$ws.Range("J45").Value = 4179
$ws.Range("K45").Value = 2482.8572
$ws.Range("L45").Value = 4179
$ws.Range("M45").Value = -2105.8572
$ws.Range("N45").Value = -4933

$ws.Range("H74").Value = 45461716
$ws.Range("I74").Value = 50007690
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 50007690
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -50006816
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 45461716
$ws.Range("I77").Value = 50007690
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 250038450
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -250034082
$ws.Range("N77").Value = -18736

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H116").Value = 2217.76
$ws.Range("I116").Value = 1349
$ws.Range("J116").Value = 4968.8335
$ws.Range("K116").Value = 1349
$ws.Range("L116").Value = 4968.8335
$ws.Range("M116").Value = 945
$ws.Range("N116").Value = -9556.833500000001

$ws.Range("H139").Value = 27238.334
$ws.Range("J139").Value = 27238.334
$ws.Range("L139").Value = 27238.334
$ws.Range("N139").Value = -37518.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2217.76
$ws.Range("I3").Value = 1349
$ws.Range("J3").Value = 4968.8335
$ws.Range("K3").Value = 1349
$ws.Range("L3").Value = 4968.8335
$ws.Range("M3").Value = -1235
$ws.Range("N3").Value = -5196.8335

$ws.Range("H86").Value = 1383.75
$ws.Range("I86").Value = 1400.2142
$ws.Range("J86").Value = 1360.7
$ws.Range("K86").Value = 1400.2142
$ws.Range("L86").Value = 1360.7
$ws.Range("M86").Value = -277.2141999999999
$ws.Range("N86").Value = -3606.7

$ws.Range("H89").Value = 1383.75
$ws.Range("I89").Value = 1400.2142
$ws.Range("J89").Value = 1360.7
$ws.Range("K89").Value = 7001.071
$ws.Range("L89").Value = 6803.5
$ws.Range("M89").Value = -1385.071
$ws.Range("N89").Value = -18035.5

$ws.Range("H105").Value = 5259
$ws.Range("I105").Value = 2062
$ws.Range("J105").Value = 8456
$ws.Range("K105").Value = 2062
$ws.Range("L105").Value = 8456
$ws.Range("M105").Value = -315
$ws.Range("N105").Value = -11950

$ws.Range("H134").Value = 6247.636
$ws.Range("I134").Value = 3539.8462
$ws.Range("J134").Value = 8007.7
$ws.Range("K134").Value = 10619.5386
$ws.Range("L134").Value = 24023.1
$ws.Range("M134").Value = -8084.5386
$ws.Range("N134").Value = -29093.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6062669.5
$ws.Range("I31").Value = 1419.4773
$ws.Range("J31").Value = 30307670
$ws.Range("K31").Value = 1419.4773
$ws.Range("L31").Value = 30307670
$ws.Range("M31").Value = -1124.4773
$ws.Range("N31").Value = -30308260

$ws.Range("H34").Value = 6062669.5
$ws.Range("I34").Value = 1419.4773
$ws.Range("J34").Value = 30307670
$ws.Range("K34").Value = 1419.4773
$ws.Range("L34").Value = 30307670
$ws.Range("M34").Value = -1217.4773
$ws.Range("N34").Value = -30308074

$ws.Range("H141").Value = 40501
$ws.Range("J141").Value = 40501
$ws.Range("L141").Value = 40501
$ws.Range("N141").Value = -50861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2214
$ws.Range("I5").Value = 510.25
$ws.Range("K5").Value = 1530.75
$ws.Range("M5").Value = -1418.75

$ws.Range("H122").Value = 3425.0825
$ws.Range("I122").Value = 343.45456
$ws.Range("J122").Value = 3819.2441
$ws.Range("K122").Value = 3091.09104
$ws.Range("L122").Value = 34373.1969
$ws.Range("M122").Value = -641.0910400000002
$ws.Range("N122").Value = -39273.1969

$ws.Range("H135").Value = 2214
$ws.Range("I135").Value = 510.25
$ws.Range("K135").Value = 4592.25
$ws.Range("M135").Value = -2057.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 508.7857
$ws.Range("I107").Value = 428.83334
$ws.Range("J107").Value = 652.7
$ws.Range("K107").Value = 428.83334
$ws.Range("L107").Value = 652.7
$ws.Range("M107").Value = 1491.16666
$ws.Range("N107").Value = -4492.7

$ws.Range("H113").Value = 2088.7778
$ws.Range("I113").Value = 1849.875
$ws.Range("K113").Value = 1849.875
$ws.Range("M113").Value = 320.125

$ws.Range("H122").Value = 359153.8
$ws.Range("I122").Value = 667892.2
$ws.Range("J122").Value = 2917.2307
$ws.Range("K122").Value = 2003676.6
$ws.Range("L122").Value = 8751.6921
$ws.Range("M122").Value = -2001226.6
$ws.Range("N122").Value = -13651.6921

$ws.Range("H126").Value = 10002210
$ws.Range("I126").Value = 15626790
$ws.Range("J126").Value = 2955.3333
$ws.Range("K126").Value = 46880370
$ws.Range("L126").Value = 8865.999899999999
$ws.Range("M126").Value = -46877900
$ws.Range("N126").Value = -13805.9999

$ws.Range("H132").Value = 1813891.9
$ws.Range("I132").Value = 3207192.5
$ws.Range("J132").Value = 2601.2
$ws.Range("K132").Value = 9621577.5
$ws.Range("L132").Value = 7803.599999999999
$ws.Range("M132").Value = -9619047.5
$ws.Range("N132").Value = -12863.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 12288.889
$ws.Range("J23").Value = 21398.8
$ws.Range("L23").Value = 21398.8
$ws.Range("N23").Value = -21858.8

$ws.Range("H55").Value = 24083.691
$ws.Range("I55").Value = 235.2
$ws.Range("K55").Value = 235.2
$ws.Range("M55").Value = -62.19999999999999

$ws.Range("H132").Value = 55470.75
$ws.Range("I132").Value = 68094.5
$ws.Range("J132").Value = 4975.75
$ws.Range("K132").Value = 204283.5
$ws.Range("L132").Value = 14927.25
$ws.Range("M132").Value = -201753.5
$ws.Range("N132").Value = -19987.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49232.25
$ws.Range("J46").Value = 49232.25
$ws.Range("L46").Value = 49232.25
$ws.Range("N46").Value = -49694.25

$ws.Range("H101").Value = 11540.167
$ws.Range("J101").Value = 11540.167
$ws.Range("L101").Value = 11540.167
$ws.Range("N101").Value = -18030.167

$ws.Range("H107").Value = 993.36365
$ws.Range("I107").Value = 785.8
$ws.Range("K107").Value = 2357.4
$ws.Range("M107").Value = -437.3999999999996

$ws.Range("H122").Value = 4745.3076
$ws.Range("I122").Value = 5526.1816
$ws.Range("J122").Value = 3389.0527
$ws.Range("K122").Value = 16578.5448
$ws.Range("L122").Value = 10167.1581
$ws.Range("M122").Value = -14128.5448
$ws.Range("N122").Value = -15067.1581

$ws.Range("H134").Value = 49232.25
$ws.Range("J134").Value = 49232.25
$ws.Range("L134").Value = 147696.75
$ws.Range("N134").Value = -152766.75

$ws.Range("H136").Value = 2123.5715
$ws.Range("I136").Value = 1220.909
$ws.Range("J136").Value = 5433.3335
$ws.Range("K136").Value = 3662.727
$ws.Range("L136").Value = 16300.0005
$ws.Range("M136").Value = -1112.727
$ws.Range("N136").Value = -21400.0005
